# Cotações atualizadas - 2025-10-07
# Adds a new row (33) with the latest daily quotes for the funds,
# continuing the existing table in Sheet1 (A1:E32 -> A1:E33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the "Data" column) so the new
# entry is appended right after it, regardless of how many rows
# already exist.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Date column: keep the same date/number format as the row above it.
$ws.Cells.Item($newRow, 1).Value = 45937
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# Quote columns, stored as text (Portuguese-style decimal comma),
# matching the existing rows in the sheet.
$ws.Cells.Item($newRow, 2).Value = "21,5982"
$ws.Cells.Item($newRow, 3).Value = "15,4244"
$ws.Cells.Item($newRow, 4).Value = "15,4503"
$ws.Cells.Item($newRow, 5).Value = "15,4503"

Write-Output "Added row $newRow with 2025-10-07 quotes"
